$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# COSY sheet: shrink the COSY-range sub-table (L:V) and the molecule-peak
# sub-table (A:I, third block) by 5 rows each, matching updated peak counts.
# ---------------------------------------------------------------------------
$cosy = $wb.Worksheets.Item("COSY")

# Copy formatting (borders) onto the new boundary cells BEFORE the rows that
# currently hold that formatting get deleted.
$cosy.Range("J2").Copy() | Out-Null
$cosy.Range("L2").PasteSpecial(-4122) | Out-Null          # xlPasteFormats

$cosy.Range("L18:V18").Copy() | Out-Null
$cosy.Range("L13:V13").PasteSpecial(-4122) | Out-Null     # new bottom border row for L:V table

$cosy.Range("U17:V17").Copy() | Out-Null
$cosy.Range("U12:V12").PasteSpecial(-4122) | Out-Null     # new bottom border row for U:V pair

$excel.CutCopyMode = 0

# Remove the 5 now-unused rows; everything below shifts up and the sheet
# dimension shrinks from A1:V35 to A1:V30.
$cosy.Rows("14:18").Delete() | Out-Null

# The sheet becomes the active tab, selection collapses to the main table.
$cosy.Activate()
$cosy.Range("A1:K13").Select() | Out-Null
